# "definitions superceded by code modules"
#
# The "Published" worksheet tracks, per definition, whether a reusable code
# module now exists for it (column H, "Code Module") and any supporting
# notes (column I, "Notes"). A handful of definitions have since been
# superseded by code modules, so their status flips from "N" to "Y" (or
# "Partial" for the driver-licence definition, which only has partial
# coverage), with notes added where relevant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Published")

# Row 99 (driver licence definition): only partially superseded by a code
# module, so note this explicitly.
$ws.Cells.Item(99, 8).Value = "Partial"
$ws.Cells.Item(99, 9).Value = "Driver licence code module now available"

# Definitions that now have a full code module available: flip "Code
# Module" from "N" to "Y".
$rowsWithCodeModule = @(3, 17, 20, 91, 92, 96, 97)
foreach ($r in $rowsWithCodeModule) {
    $ws.Cells.Item($r, 8).Value = "Y"
}

# Rows 96 and 97 also get an explanatory note.
$ws.Cells.Item(96, 9).Value = "Code module now available"
$ws.Cells.Item(97, 9).Value = "Code module now available"
